# "Merged PR 138: updated excel with Test1 env information"
#
# Adds a new "Env" column to the "Server List" and "Security Groups" sheets
# (marking all of the pre-existing rows "Dev"), appends the matching "Test1"
# environment rows to each, re-applies AutoFilter/sort, and leaves
# "Security Groups" as the active sheet/tab.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Server List": insert an "Env" column at A (existing Type/Name shift
# to B/C), tag the pre-existing rows "Dev", type the new Test1 rows in the
# order they were entered (Job*, ETL, Tidal, Web*), then sort the table by
# Type so the final row order is ETL / Job*6 / Tidal / Web*6.
# ---------------------------------------------------------------------------
$wsServers = $wb.Worksheets.Item("Server List")
$wsServers.Columns.Item(1).Insert()

$wsServers.Range("A1").Value = "Env"
$wsServers.Range("A2:A9").Value = "Dev"

$serverRows = @(
    @("Test1", "Job",   "IRMATest1Job01"),
    @("Test1", "Job",   "IRMATest1Job02"),
    @("Test1", "Job",   "IRMATest1Job03"),
    @("Test1", "Job",   "IRMATest1Job04"),
    @("Test1", "Job",   "IRMATest1Job05"),
    @("Test1", "Job",   "IRMATest1Job06"),
    @("Test1", "ETL",   "IRMATest1ETL01"),
    @("Test1", "Tidal", "IRMATest1Tidal01"),
    @("Test1", "Web",   "IRMATest1Web01"),
    @("Test1", "Web",   "IRMATest1Web02"),
    @("Test1", "Web",   "IRMATest1Web03"),
    @("Test1", "Web",   "IRMATest1Web04"),
    @("Test1", "Web",   "IRMATest1Web05"),
    @("Test1", "Web",   "IRMATest1Web06")
)

$r = 10
foreach ($row in $serverRows) {
    $wsServers.Cells.Item($r, 1).Value = $row[0]
    $wsServers.Cells.Item($r, 2).Value = $row[1]
    $wsServers.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Sort the Test1 block into alphabetical Type order (ETL, Job, Tidal, Web).
$sortRange = $wsServers.Range("A10:C23")
$wsServers.Sort.SortFields.Clear() | Out-Null
$wsServers.Sort.SortFields.Add($wsServers.Range("B10:B23")) | Out-Null
$wsServers.Sort.SetRange($sortRange) | Out-Null
$wsServers.Sort.Header = 2
$wsServers.Sort.Apply() | Out-Null
$wsServers.Sort.SortFields.Clear() | Out-Null

$wsServers.Range("A1:C1").AutoFilter() | Out-Null
$fdServers = $wsServers.Names.Add("_xlnm._FilterDatabase", "='Server List'!`$A`$1:`$C`$1")
$fdServers.Visible = $false

$wsServers.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Sheet "Security Groups": insert an "Env" column at A (existing columns
# shift right by one), tag the pre-existing rows "Dev", then append the
# Test1 rows (already in final row order).
# ---------------------------------------------------------------------------
$wsGroups = $wb.Worksheets.Item("Security Groups")
$wsGroups.Columns.Item(1).Insert()

$wsGroups.Range("A1").Value = "Env"
$wsGroups.Range("A2:A14").Value = "Dev"

$groupRows = @(
    @("Test1", "Administrators", "WFM\IRMA.Deploy.NonProd",  "Job Web"),
    @("Test1", "Administrators", "WFM\Icon.Deploy.NonProd",  "Job Web"),
    @("Test1", "Administrators", "WFM\Mammoth.Deploy.NProd", "Job Web"),
    @("Test1", "Administrators", "WFM\IconInterfaceUserTes", "Job Web"),
    @("Test1", "Administrators", "wfm\MammothTest",          "Job Web"),
    @("Test1", "Administrators", "wfm\IconWebTest",          "Web"),
    @("Test1", "IIS_IUSRS",      "wfm\MammothTest",          "Web"),
    @("Test1", "IIS_IUSRS",      "wfm\IconWebTest",          "Web"),
    @("Test1", "IIS_IUSRS",      "Authenticated Users",      "Web")
)

$r = 15
foreach ($row in $groupRows) {
    $wsGroups.Cells.Item($r, 1).Value = $row[0]
    $wsGroups.Cells.Item($r, 2).Value = $row[1]
    $wsGroups.Cells.Item($r, 3).Value = $row[2]
    $wsGroups.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$wsGroups.Range("A1:D1").AutoFilter() | Out-Null
$fdGroups = $wsGroups.Names.Add("_xlnm._FilterDatabase", "='Security Groups'!`$A`$1:`$D`$1")
$fdGroups.Visible = $false

# ---------------------------------------------------------------------------
# Leave selections/active tab the way the author left them: "Server List"
# parked at G24 and "Security Groups" active/selected at J11.
# ---------------------------------------------------------------------------
$wsServers.Range("G24").Select() | Out-Null
$wsGroups.Select() | Out-Null
$wsGroups.Range("J11").Select() | Out-Null
